# Apply transcript renaming edits to the DataSheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: "Student Tag" -> "S Tag"
$ws.Range("G1").Value = "S Tag"

# Column D rows where "RBD" -> "T"
$rbdRows = @(6,7,8,12,13,15,16,17,18,19,20,21,22,27,28,30,31,33,35,37,39,40,42,44)
foreach ($r in $rbdRows) {
    $ws.Cells.Item($r, 4).Value = "T"
}

# Column D rows where "Student" -> "S"
$studentRows = @(9,41)
foreach ($r in $studentRows) {
    $ws.Cells.Item($r, 4).Value = "S"
}

# Column G rows where "2 - relating to another student" -> "2 - relating to another S"
$ws.Range("G9").Value = "2 - relating to another S"
$ws.Range("G11").Value = "2 - relating to another S"
